# "add slide code review"
#
# 1. Finish the "Code Review?" slide title -> "Code Review" (drop the
#    trailing question mark and tidy it up into a single run).
# 2. Reorder the two closing slides so "Take away" now comes right
#    before "Code Review" (they were swapped from their previous order).

$p = $ppt.ActivePresentation

# --- 1. Clean up the "Code Review?" slide title -----------------------
$codeReviewSlide = $p.Slides.Item(10)
$title = $codeReviewSlide.Shapes.Item(1)
$tr = $title.TextFrame.TextRange
$tr.Text = ""
$tr.Text = "Code Review"

# --- 2. Move the "Take away" slide so it precedes "Code Review" -------
$takeAwaySlide = $p.Slides.Item(11)
$takeAwaySlide.MoveTo(10)
